# Added CC tip - mental mappings
#
# 1) Every "Date Placeholder" (master + all slide layouts) has its cached
#    datetimeFigureOut field text refreshed from 11/26/2021 -> 1/12/2022
#    (PowerPoint re-caches this automatic field on every save).
# 2) Slide 2's cover title is updated from "Tip #7 / Keep parameters in a
#    consistent order" to "Tip #8 / avoid mental mappings".

$p = $ppt.ActivePresentation

function Update-DateShapes($shapes, $newDateText) {
    $n = $shapes.Count
    for ($i = 1; $i -le $n; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder) {
            $tr = $shp.TextFrame.TextRange
            $curLen = $tr.Text.Length
            if ($curLen -gt 0) {
                $whole = $tr.Characters(1, $curLen)
                $whole.Text = $newDateText
            } else {
                $tr.Text = $newDateText
            }
        }
    }
}

# --- 1. Refresh the cached "today" field on the slide master ---
$master = $p.SlideMaster
Update-DateShapes $master.Shapes "1/12/2022"

# --- 1b. ... and on every slide layout that hangs off the master ---
$layouts = $master.CustomLayouts
$layoutCount = $layouts.Count
for ($L = 1; $L -le $layoutCount; $L++) {
    $lyt = $layouts.Item($L)
    Update-DateShapes $lyt.Shapes "1/12/2022"
}

# --- 2. Update the "Tip #7" cover slide text to "Tip #8" ---
$slide = $p.Slides.Item(2)
$title = $slide.Shapes.Item(2)
$titleRange = $title.TextFrame.TextRange

$tipHeading = $titleRange.Find("Tip #7")
if ($tipHeading -ne $null) {
    $tipHeading.Text = "Tip #8"
}

$tipBody = $titleRange.Find("Keep parameters in a consistent order")
if ($tipBody -ne $null) {
    $tipBody.Text = "avoid mental mappings"
}
